# Atualização automática de preços de eletricidade
# Update the daily hourly spot-price row (row 2) on the Spot_PT sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spot_PT")

# Day serial (keeps existing date number-format on the cell)
$ws.Range("A2").Value = 46004

# Hourly prices 0h-1h ... 23h-24h (columns B..Y)
$ws.Range("B2").Value = 90.73
$ws.Range("C2").Value = 81.09
$ws.Range("D2").Value = 74.39
$ws.Range("E2").Value = 68.95
$ws.Range("F2").Value = 65.73999999999999
$ws.Range("G2").Value = 67.12
$ws.Range("H2").Value = 70.58
$ws.Range("I2").Value = 78.31
$ws.Range("J2").Value = 84.95999999999999
$ws.Range("K2").Value = 80.79000000000001
$ws.Range("L2").Value = 72.44
$ws.Range("M2").Value = 67.64
$ws.Range("N2").Value = 65.52
$ws.Range("O2").Value = 59.31
$ws.Range("P2").Value = 59.98
$ws.Range("Q2").Value = 66.42
$ws.Range("R2").Value = 76.62
$ws.Range("S2").Value = 91.91
$ws.Range("T2").Value = 99.41
$ws.Range("U2").Value = 96.28
$ws.Range("V2").Value = 97.90000000000001
$ws.Range("W2").Value = 98.62
$ws.Range("X2").Value = 98.23
$ws.Range("Y2").Value = 87.81

# Daily average price
$ws.Range("Z2").Value = 79.2

# Max 4h slot window + price
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 95.64

# First 2h slot (within the max 4h window) + price
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 98.26000000000001

# Second 2h slot (within the max 4h window) + price -- unchanged window, price updated
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 97.84

# Cheapest slot window
$ws.Range("AG2").Value = "2h-16h"
